# Auto-generated script applying scheduled market-data refresh values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 94.46154
$ws.Range("I8").Value = 67.90000000000001
$ws.Range("K8").Value = 203.7
$ws.Range("M8").Value = -64.70000000000002
$ws.Range("H18").Value = 1274.091
$ws.Range("I18").Value = 1246.5
$ws.Range("J18").Value = 1550
$ws.Range("K18").Value = 1246.5
$ws.Range("L18").Value = 1550
$ws.Range("M18").Value = -962.5
$ws.Range("N18").Value = -2118
$ws.Range("H28").Value = 1149.75
$ws.Range("I28").Value = 690.75
$ws.Range("K28").Value = 690.75
$ws.Range("M28").Value = -205.75
$ws.Range("H33").Value = 453.94116
$ws.Range("I33").Value = 249.7931
$ws.Range("K33").Value = 249.7931
$ws.Range("M33").Value = -20.79310000000001
$ws.Range("H64").Value = 5832
$ws.Range("J64").Value = 5623.5
$ws.Range("L64").Value = 5623.5
$ws.Range("N64").Value = -6119.5
$ws.Range("H67").Value = 5832
$ws.Range("J67").Value = 5623.5
$ws.Range("L67").Value = 5623.5
$ws.Range("N67").Value = -7339.5
$ws.Range("H87").Value = 123177.73
$ws.Range("J87").Value = 123177.73
$ws.Range("L87").Value = 123177.73
$ws.Range("N87").Value = -125673.73
$ws.Range("H90").Value = 123177.73
$ws.Range("J90").Value = 123177.73
$ws.Range("L90").Value = 369533.19
$ws.Range("N90").Value = -382013.19
$ws.Range("H107").Value = 41864.867
$ws.Range("I107").Value = 34665.332
$ws.Range("J107").Value = 52664.168
$ws.Range("K107").Value = 34665.332
$ws.Range("L107").Value = 52664.168
$ws.Range("M107").Value = -32745.332
$ws.Range("N107").Value = -56504.168
$ws.Range("H111").Value = 3736.476
$ws.Range("I111").Value = 1803.5834
$ws.Range("J111").Value = 6313.6665
$ws.Range("K111").Value = 5410.7502
$ws.Range("L111").Value = 18940.9995
$ws.Range("M111").Value = -2343.7502
$ws.Range("N111").Value = -25074.9995
$ws.Range("H112").Value = 3311.5881
$ws.Range("I112").Value = 944.25
$ws.Range("J112").Value = 4040
$ws.Range("K112").Value = 2832.75
$ws.Range("L112").Value = 12120
$ws.Range("M112").Value = -1724.75
$ws.Range("N112").Value = -14336
$ws.Range("H113").Value = 1223.75
$ws.Range("I113").Value = 1231.6666
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1231.6666
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 2022.3334
$ws.Range("N113").Value = -7708
$ws.Range("H118").Value = 550.5294
$ws.Range("I118").Value = 588.4167
$ws.Range("K118").Value = 1765.2501
$ws.Range("M118").Value = -108.2501
$ws.Range("H138").Value = 6438.9414
$ws.Range("J138").Value = 7525.3706
$ws.Range("L138").Value = 22576.1118
$ws.Range("N138").Value = -32856.1118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7308.1924
$ws.Range("I32").Value = 756.18604
$ws.Range("K32").Value = 756.18604
$ws.Range("M32").Value = -469.18604
$ws.Range("H97").Value = 1812.2858
$ws.Range("J97").Value = 2762.111
$ws.Range("L97").Value = 2762.111
$ws.Range("N97").Value = -3754.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 26569.334
$ws.Range("J60").Value = 29499.5
$ws.Range("L60").Value = 29499.5
$ws.Range("N60").Value = -30697.5
$ws.Range("H99").Value = 3027.8462
$ws.Range("I99").Value = 2851.4546
$ws.Range("J99").Value = 3998
$ws.Range("K99").Value = 2851.4546
$ws.Range("L99").Value = 3998
$ws.Range("M99").Value = -1353.4546
$ws.Range("N99").Value = -6994
$ws.Range("H105").Value = 4082.7058
$ws.Range("I105").Value = 4138.2
$ws.Range("K105").Value = 4138.2
$ws.Range("M105").Value = -2391.2
$ws.Range("H107").Value = 2264.9412
$ws.Range("I107").Value = 1938.8462
$ws.Range("J107").Value = 3324.75
$ws.Range("K107").Value = 1938.8462
$ws.Range("L107").Value = 3324.75
$ws.Range("M107").Value = -18.84619999999995
$ws.Range("N107").Value = -7164.75
$ws.Range("H134").Value = 2424.4666
$ws.Range("I134").Value = 2240.5
$ws.Range("K134").Value = 6721.5
$ws.Range("M134").Value = -4186.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 59809
$ws.Range("J16").Value = 48881.668
$ws.Range("L16").Value = 48881.668
$ws.Range("N16").Value = -49455.668
$ws.Range("H42").Value = 11499.5
$ws.Range("I42").Value = 11499.5
$ws.Range("K42").Value = 11499.5
$ws.Range("M42").Value = -10906.5
$ws.Range("H58").Value = 3789.8
$ws.Range("I58").Value = 3679.375
$ws.Range("K58").Value = 3679.375
$ws.Range("M58").Value = -3476.375
$ws.Range("H113").Value = 59809
$ws.Range("J113").Value = 48881.668
$ws.Range("L113").Value = 48881.668
$ws.Range("N113").Value = -53221.668
$ws.Range("H136").Value = 3789.8
$ws.Range("I136").Value = 3679.375
$ws.Range("K136").Value = 11038.125
$ws.Range("M136").Value = -8488.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3757.25
$ws.Range("I3").Value = 3757.25
$ws.Range("K3").Value = 11271.75
$ws.Range("M3").Value = -11159.75
$ws.Range("H131").Value = 1852.3334
$ws.Range("I131").Value = 1008.6
$ws.Range("J131").Value = 2176.8462
$ws.Range("K131").Value = 3025.8
$ws.Range("L131").Value = 6530.5386
$ws.Range("M131").Value = 2014.2
$ws.Range("N131").Value = -16610.5386
$ws.Range("H139").Value = 3689.25
$ws.Range("I139").Value = 2800.1875
$ws.Range("J139").Value = 7245.5
$ws.Range("K139").Value = 8400.5625
$ws.Range("L139").Value = 21736.5
$ws.Range("M139").Value = -3260.5625
$ws.Range("N139").Value = -32016.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 700000
$ws.Range("J4").Value = 400000
$ws.Range("L4").Value = 400000
$ws.Range("N4").Value = -400224
$ws.Range("H132").Value = 2481.5557
$ws.Range("I132").Value = 2301
$ws.Range("J132").Value = 2765.2856
$ws.Range("K132").Value = 6903
$ws.Range("L132").Value = 8295.856800000001
$ws.Range("M132").Value = -4373
$ws.Range("N132").Value = -13355.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10592.714
$ws.Range("I40").Value = 2310.6667
$ws.Range("J40").Value = 25500.4
$ws.Range("K40").Value = 2310.6667
$ws.Range("L40").Value = 25500.4
$ws.Range("M40").Value = -2174.6667
$ws.Range("N40").Value = -25772.4
$ws.Range("H68").Value = 3815.3333
$ws.Range("I68").Value = 3815.3333
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3815.3333
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -3066.3333
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3815.3333
$ws.Range("I71").Value = 3815.3333
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 19076.6665
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -15332.6665
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 66040.09
$ws.Range("I132").Value = 79937.89
$ws.Range("K132").Value = 239813.67
$ws.Range("M132").Value = -237283.67
$ws.Range("H136").Value = 6812.7
$ws.Range("I136").Value = 6826.3335
$ws.Range("J136").Value = 6690
$ws.Range("K136").Value = 20479.0005
$ws.Range("L136").Value = 20070
$ws.Range("M136").Value = -17929.0005
$ws.Range("N136").Value = -25170
$ws.Range("H137").Value = 73331.664
$ws.Range("J137").Value = 73331.664
$ws.Range("L137").Value = 73331.664
$ws.Range("N137").Value = -83531.664
$ws.Range("H138").Value = 112899.75
$ws.Range("J138").Value = 112899.75
$ws.Range("L138").Value = 112899.75
$ws.Range("N138").Value = -123179.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2559.2
$ws.Range("I107").Value = 2455.375
$ws.Range("K107").Value = 7366.125
$ws.Range("M107").Value = -5446.125
$ws.Range("H132").Value = 2837.1667
$ws.Range("I132").Value = 1721.1428
$ws.Range("K132").Value = 5163.428400000001
$ws.Range("M132").Value = -2633.428400000001
